# Update "想去人数" (interest count, column F) values for several rows
# across the 展览 (sheet 1), 演出 (sheet 2) and 全部类型 (sheet 4) sheets.
# This mirrors an automated re-scrape of source data (see commit message).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 610
$ws1.Range("F6").Value  = 15192
$ws1.Range("F10").Value = 15291
$ws1.Range("F12").Value = 8830
$ws1.Range("F13").Value = 352
$ws1.Range("F16").Value = 186
$ws1.Range("F20").Value = 36
$ws1.Range("F24").Value = 56
$ws1.Range("F29").Value = 31
$ws1.Range("F32").Value = 38
$ws1.Range("F35").Value = 285
$ws1.Range("F36").Value = 434
$ws1.Range("F38").Value = 5419

# --- Sheet "演出" --------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 63

# --- Sheet "全部类型" ----------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 610
$ws4.Range("F6").Value  = 15192
$ws4.Range("F10").Value = 15291
$ws4.Range("F12").Value = 8830
$ws4.Range("F13").Value = 352
$ws4.Range("F17").Value = 186
$ws4.Range("F21").Value = 36
$ws4.Range("F25").Value = 56
$ws4.Range("F30").Value = 31
$ws4.Range("F32").Value = 63
$ws4.Range("F35").Value = 38
$ws4.Range("F38").Value = 285
$ws4.Range("F39").Value = 434
$ws4.Range("F41").Value = 5419
